$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1107.4375
$ws.Range("I2").Value = 144.2
$ws.Range("J2").Value = 2712.8333
$ws.Range("K2").Value = 144.2
$ws.Range("L2").Value = 2712.8333
$ws.Range("M2").Value = -31.19999999999999
$ws.Range("N2").Value = -2938.8333
$ws.Range("H4").Value = 227
$ws.Range("I4").Value = 227
$ws.Range("K4").Value = 227
$ws.Range("M4").Value = -113
$ws.Range("H6").Value = 148.6
$ws.Range("I6").Value = 148.6
$ws.Range("K6").Value = 445.8
$ws.Range("M6").Value = -333.8
$ws.Range("H28").Value = 9465.714
$ws.Range("I28").Value = 7861.643
$ws.Range("J28").Value = 12673.857
$ws.Range("K28").Value = 7861.643
$ws.Range("L28").Value = 12673.857
$ws.Range("M28").Value = -7376.643
$ws.Range("N28").Value = -13643.857
$ws.Range("H88").Value = 11999
$ws.Range("J88").Value = 11999
$ws.Range("L88").Value = 11999
$ws.Range("N88").Value = -12811
$ws.Range("H91").Value = 11999
$ws.Range("J91").Value = 11999
$ws.Range("L91").Value = 11999
$ws.Range("N91").Value = -14807
$ws.Range("H132").Value = 465062.72
$ws.Range("I132").Value = 539656.6
$ws.Range("K132").Value = 1618969.8
$ws.Range("M132").Value = -1616439.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 374.4
$ws.Range("J4").Value = 550
$ws.Range("L4").Value = 550
$ws.Range("N4").Value = -782
$ws.Range("H5").Value = 625.5
$ws.Range("I5").Value = 625.5
$ws.Range("K5").Value = 625.5
$ws.Range("M5").Value = -513.5
$ws.Range("H32").Value = 1357226.4
$ws.Range("I32").Value = 869.0625
$ws.Range("K32").Value = 869.0625
$ws.Range("M32").Value = -582.0625
$ws.Range("H43").Value = 73249.75
$ws.Range("J43").Value = 73249.75
$ws.Range("L43").Value = 73249.75
$ws.Range("N43").Value = -73875.75
$ws.Range("H61").Value = 5081.423
$ws.Range("I61").Value = 5129.3335
$ws.Range("J61").Value = 4506.5
$ws.Range("K61").Value = 5129.3335
$ws.Range("L61").Value = 4506.5
$ws.Range("M61").Value = -4917.3335
$ws.Range("N61").Value = -4930.5
$ws.Range("H74").Value = 4532.1577
$ws.Range("I74").Value = 5131.1665
$ws.Range("K74").Value = 5131.1665
$ws.Range("M74").Value = -4257.1665
$ws.Range("H77").Value = 4532.1577
$ws.Range("I77").Value = 5131.1665
$ws.Range("K77").Value = 25655.8325
$ws.Range("M77").Value = -21287.8325
$ws.Range("H136").Value = 5081.423
$ws.Range("I136").Value = 5129.3335
$ws.Range("J136").Value = 4506.5
$ws.Range("K136").Value = 15388.0005
$ws.Range("L136").Value = 13519.5
$ws.Range("M136").Value = -12838.0005
$ws.Range("N136").Value = -18619.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 625.5
$ws.Range("I4").Value = 625.5
$ws.Range("K4").Value = 625.5
$ws.Range("M4").Value = -510.5
$ws.Range("H105").Value = 3131.5676
$ws.Range("J105").Value = 3418.88
$ws.Range("L105").Value = 3418.88
$ws.Range("N105").Value = -6912.88

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 206.63158
$ws.Range("I7").Value = 152.90909
$ws.Range("K7").Value = 152.90909
$ws.Range("M7").Value = -39.90908999999999
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H62").Value = 11910.2
$ws.Range("I62").Value = 13888
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 13888
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -13264
$ws.Range("N62").Value = -5247
$ws.Range("H65").Value = 11910.2
$ws.Range("I65").Value = 13888
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 69440
$ws.Range("L65").Value = 19995
$ws.Range("M65").Value = -66320
$ws.Range("N65").Value = -26235
$ws.Range("H134").Value = 100010410
$ws.Range("I134").Value = 111119370
$ws.Range("J134").Value = 29750
$ws.Range("K134").Value = 333358110
$ws.Range("L134").Value = 89250
$ws.Range("M134").Value = -333355575
$ws.Range("N134").Value = -94320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 313442.2
$ws.Range("I5").Value = 317.66666
$ws.Range("J5").Value = 385701.7
$ws.Range("K5").Value = 952.9999799999999
$ws.Range("L5").Value = 1157105.1
$ws.Range("M5").Value = -840.9999799999999
$ws.Range("N5").Value = -1157329.1
$ws.Range("H94").Value = 5626.6665
$ws.Range("J94").Value = 6399.4
$ws.Range("L94").Value = 19198.2
$ws.Range("N94").Value = -20550.2
$ws.Range("H132").Value = 1998.3846
$ws.Range("I132").Value = 1005.625
$ws.Range("J132").Value = 3586.8
$ws.Range("K132").Value = 9050.625
$ws.Range("L132").Value = 32281.2
$ws.Range("M132").Value = -6520.625
$ws.Range("N132").Value = -37341.2
$ws.Range("H135").Value = 313442.2
$ws.Range("I135").Value = 317.66666
$ws.Range("J135").Value = 385701.7
$ws.Range("K135").Value = 2858.99994
$ws.Range("L135").Value = 3471315.3
$ws.Range("M135").Value = -323.9999399999997
$ws.Range("N135").Value = -3476385.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 42714.285
$ws.Range("J117").Value = 42714.285
$ws.Range("L117").Value = 42714.285
$ws.Range("N117").Value = -49598.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5480.591
$ws.Range("I7").Value = 4795
$ws.Range("K7").Value = 4795
$ws.Range("M7").Value = -4683
$ws.Range("H100").Value = 1500
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H126").Value = 5480.591
$ws.Range("I126").Value = 4795
$ws.Range("K126").Value = 14385
$ws.Range("M126").Value = -11915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 12515368
$ws.Range("I13").Value = 33334300
$ws.Range("J13").Value = 24008.8
$ws.Range("K13").Value = 33334300
$ws.Range("L13").Value = 24008.8
$ws.Range("M13").Value = -33334160
$ws.Range("N13").Value = -24288.8
$ws.Range("H41").Value = 13418.777
$ws.Range("J41").Value = 13177.625
$ws.Range("L41").Value = 13177.625
$ws.Range("N41").Value = -13957.625
